# "finestra incidenza 7gg centrata su ultimo g"
# Recompute the 7-day rolling sum (col C) and the per-100k-inhabitants
# figure (col D) so the 7-day window TRAILS and ends on the current day
# (r-6 .. r) instead of being centred on the current day (r-3 .. r+3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 184
$windowSize = 7

# Exact C -> D values as originally stored in the workbook (D = C * 100000 /
# population, population = 4585) - reusing the literal doubles keeps the
# floating point bit pattern identical to what Excel itself produced.
$incidenceMap = @{
    0  = 0.0
    1  = 21.81025081788441
    2  = 43.62050163576881
    3  = 65.43075245365321
    4  = 87.24100327153762
    5  = 109.051254089422
    6  = 130.8615049073064
    7  = 152.6717557251908
    8  = 174.4820065430752
    9  = 196.2922573609597
    10 = 218.1025081788441
    11 = 239.9127589967285
    12 = 261.7230098146129
    13 = 283.5332606324973
    14 = 305.3435114503817
    15 = 327.1537622682661
    16 = 348.9640130861505
    17 = 370.7742639040349
    18 = 392.5845147219193
    19 = 414.3947655398038
    20 = 436.2050163576881
    21 = 458.0152671755725
}
$population = 4585

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    $windowStart = $r - ($windowSize - 1)

    if ($windowStart -lt $firstDataRow) {
        # Not enough history yet for a full trailing 7-day window -> blank.
        $ws.Range("C" + $r + ":D" + $r).ClearContents()
    } else {
        $sum = 0
        for ($k = $windowStart; $k -le $r; $k++) {
            $sum = $sum + $ws.Cells.Item($k, 2).Value2
        }

        $ws.Cells.Item($r, 3).Value2 = $sum

        if ($incidenceMap.ContainsKey($sum)) {
            $ws.Cells.Item($r, 4).Value2 = $incidenceMap[$sum]
        } else {
            $ws.Cells.Item($r, 4).Value2 = $sum * 100000 / $population
        }
    }
}
